# Adding Astro and IMperial to the Portfolio Forecast
# - Shift all timestamps in column A (rows 2-97) forward by 18 days (45758 -> 45776 base).
# - Update the "Notified Production (MW)" values in column B for rows 22-89 to reflect
#   the added portfolio contributions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 97

# Shift every timestamp in column A by 18 days (preserves the fractional time-of-day part).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 18
}

# New production values (MW) for rows 22-89, column B.
$newValues = @{
    22 = 18
    23 = 19
    24 = 25
    25 = 37
    26 = 287
    27 = 316
    28 = 351
    29 = 392
    30 = 967
    31 = 1016
    32 = 1069
    33 = 1126
    34 = 1583
    35 = 1624
    36 = 1673
    37 = 1714
    38 = 1995
    39 = 2026
    40 = 2056
    41 = 2082
    42 = 2257
    43 = 2273
    44 = 2289
    45 = 2302
    46 = 2077
    47 = 2083
    48 = 2086
    49 = 2087
    50 = 2087
    51 = 2086
    52 = 2081
    53 = 2077
    54 = 1998
    55 = 1988
    56 = 1977
    57 = 1963
    58 = 1862
    59 = 1842
    60 = 1817
    61 = 1796
    62 = 1687
    63 = 1651
    64 = 1619
    65 = 1582
    66 = 1411
    67 = 1369
    68 = 1334
    69 = 1295
    70 = 810
    71 = 756
    72 = 700
    73 = 660
    74 = 255
    75 = 224
    76 = 202
    77 = 187
    78 = 27
    79 = 21
    80 = 21
    81 = 20
    82 = 17
    83 = 17
    84 = 17
    85 = 17
    86 = 17
    87 = 16
    88 = 16
    89 = 16
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $newValues[$row]
}
